$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.741.46"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "2.526.12"
$ws.Range("E3").Value = "  -3.61%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'585.09"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'171.31"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").Value = "2.524.71"
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").Value = "2.981.69"
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "66.516.57"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "2.528.94"
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("D19").Value = "'7.86"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("E20").Value = "  -5.20%  "
$ws.Range("D21").Value = "'348.51"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'70.13"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -3.90%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D30").Value = "0.0₃0982"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").Value = "'527.83"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D33").Value = "'1.32"
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "'157.09"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").Value = "'5.10"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'2.51"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'149.66"
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'0.560"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'3.69"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "'1.72"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0271"
$ws.Range("E50").Value = "  -9.88%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0758"
$ws.Range("E51").Value = "  -1.76%  "

Write-Host "Applied cryptos update"